# The workbook previously had a "row 7" entry duplicating the "helix jump"
# keyword/appID pair, and row 6 carried two extra hyperlinked cells
# ("Explore" / "Show (4)") pointing at asodesk.com. This edit removes the
# hyperlinks and their cell text, clears the now-empty trailing row, and
# nudges row 6's height down slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that lived on C6 and D6 (along with their backing
# relationships) before touching the cells/row that held them.
$ws.Hyperlinks.Delete()

# Clear the "Explore" / "Show (4)" text out of C6:D6 but keep their
# formatting (style index) intact.
$ws.Range("C6:D6").ClearContents()

# Drop the now-redundant row 7 entirely, shifting the dimension/used range
# back down to A1:D6.
$ws.Rows("7:7").Select()
$ws.Rows("7:7").Delete()

# Slightly shrink row 6's height.
$ws.Rows("6:6").RowHeight = 23.85
